# Auto-generated edit script: update cached Leve profit values (scheduled price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 63404.25
$ws.Range("I2").Value = 471.6
$ws.Range("K2").Value = 471.6
$ws.Range("M2").Value = -358.6
$ws.Range("H32").Value = 5280.9375
$ws.Range("I32").Value = 6645.4287
$ws.Range("J32").Value = 4219.6665
$ws.Range("K32").Value = 6645.4287
$ws.Range("L32").Value = 4219.6665
$ws.Range("M32").Value = -6319.4287
$ws.Range("N32").Value = -4871.6665
$ws.Range("H69").Value = 13843.429
$ws.Range("I69").Value = 13600.8
$ws.Range("K69").Value = 40802.39999999999
$ws.Range("M69").Value = -39928.39999999999
$ws.Range("H72").Value = 13843.429
$ws.Range("I72").Value = 13600.8
$ws.Range("K72").Value = 122407.2
$ws.Range("M72").Value = -118039.2
$ws.Range("H74").Value = 18492.309
$ws.Range("I74").Value = 4133.3335
$ws.Range("K74").Value = 4133.3335
$ws.Range("M74").Value = -3197.3335
$ws.Range("H77").Value = 18492.309
$ws.Range("I77").Value = 4133.3335
$ws.Range("K77").Value = 20666.6675
$ws.Range("M77").Value = -15986.6675
$ws.Range("H132").Value = 2610.9678
$ws.Range("I132").Value = 2636.4285
$ws.Range("K132").Value = 7909.2855
$ws.Range("M132").Value = -5379.2855
$ws.Range("H138").Value = 2606.5264
$ws.Range("J138").Value = 3421.6365
$ws.Range("L138").Value = 10264.9095
$ws.Range("N138").Value = -20544.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2888.6667
$ws.Range("I32").Value = 2648.4565
$ws.Range("K32").Value = 2648.4565
$ws.Range("M32").Value = -2361.4565
$ws.Range("H63").Value = 1766.5
$ws.Range("I63").Value = 1766.5
$ws.Range("K63").Value = 1766.5
$ws.Range("M63").Value = -1080.5
$ws.Range("H66").Value = 1766.5
$ws.Range("I66").Value = 1766.5
$ws.Range("K66").Value = 8832.5
$ws.Range("M66").Value = -5400.5
$ws.Range("H88").Value = 1515.8572
$ws.Range("I88").Value = 1628.3334
$ws.Range("J88").Value = 1431.5
$ws.Range("K88").Value = 1628.3334
$ws.Range("L88").Value = 1431.5
$ws.Range("M88").Value = -1222.3334
$ws.Range("N88").Value = -2243.5
$ws.Range("H91").Value = 1515.8572
$ws.Range("I91").Value = 1628.3334
$ws.Range("J91").Value = 1431.5
$ws.Range("K91").Value = 1628.3334
$ws.Range("L91").Value = 1431.5
$ws.Range("M91").Value = -224.3334
$ws.Range("N91").Value = -4239.5
$ws.Range("H122").Value = 1163.3704
$ws.Range("I122").Value = 1070.909
$ws.Range("K122").Value = 3212.727
$ws.Range("M122").Value = -762.7270000000003
$ws.Range("H132").Value = 2914.2
$ws.Range("I132").Value = 2552.0476
$ws.Range("K132").Value = 7656.1428
$ws.Range("M132").Value = -5126.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1982
$ws.Range("I86").Value = 2369
$ws.Range("J86").Value = 1552
$ws.Range("K86").Value = 2369
$ws.Range("L86").Value = 1552
$ws.Range("M86").Value = -1246
$ws.Range("N86").Value = -3798
$ws.Range("H89").Value = 1982
$ws.Range("I89").Value = 2369
$ws.Range("J89").Value = 1552
$ws.Range("K89").Value = 11845
$ws.Range("L89").Value = 7760
$ws.Range("M89").Value = -6229
$ws.Range("N89").Value = -18992
$ws.Range("H105").Value = 29413140
$ws.Range("I105").Value = 41667948
$ws.Range("J105").Value = 1603.8
$ws.Range("K105").Value = 41667948
$ws.Range("L105").Value = 1603.8
$ws.Range("M105").Value = -41666201
$ws.Range("N105").Value = -5097.8
$ws.Range("H134").Value = 8871.625
$ws.Range("I134").Value = 8871.625
$ws.Range("K134").Value = 26614.875
$ws.Range("M134").Value = -24079.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2144.8696
$ws.Range("I31").Value = 1036.8889
$ws.Range("J31").Value = 2857.1428
$ws.Range("K31").Value = 1036.8889
$ws.Range("L31").Value = 2857.1428
$ws.Range("M31").Value = -741.8888999999999
$ws.Range("N31").Value = -3447.1428
$ws.Range("H34").Value = 2144.8696
$ws.Range("I34").Value = 1036.8889
$ws.Range("J34").Value = 2857.1428
$ws.Range("K34").Value = 1036.8889
$ws.Range("L34").Value = 2857.1428
$ws.Range("M34").Value = -834.8888999999999
$ws.Range("N34").Value = -3261.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5709.2354
$ws.Range("I70").Value = 5532.6
$ws.Range("K70").Value = 5532.6
$ws.Range("M70").Value = -5262.6
$ws.Range("H73").Value = 5709.2354
$ws.Range("I73").Value = 5532.6
$ws.Range("K73").Value = 5532.6
$ws.Range("M73").Value = -4596.6
$ws.Range("H80").Value = 1597
$ws.Range("I80").Value = 1597
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1597
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -599
$ws.Range("H83").Value = 1597
$ws.Range("I83").Value = 1597
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7985
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2993
$ws.Range("H122").Value = 91996.28999999999
$ws.Range("I122").Value = 127959.81
$ws.Range("J122").Value = 20069.25
$ws.Range("K122").Value = 383879.43
$ws.Range("L122").Value = 60207.75
$ws.Range("M122").Value = -381429.43
$ws.Range("N122").Value = -65107.75
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4170.353
$ws.Range("I16").Value = 4146.273
$ws.Range("K16").Value = 4146.273
$ws.Range("M16").Value = -3976.273
$ws.Range("H22").Value = 4614.0835
$ws.Range("J22").Value = 5432.875
$ws.Range("L22").Value = 5432.875
$ws.Range("N22").Value = -6022.875
$ws.Range("H27").Value = 4614.0835
$ws.Range("J27").Value = 5432.875
$ws.Range("L27").Value = 5432.875
$ws.Range("N27").Value = -5646.875
$ws.Range("H68").Value = 4342.857
$ws.Range("I68").Value = 2174.625
$ws.Range("K68").Value = 2174.625
$ws.Range("M68").Value = -1425.625
$ws.Range("H71").Value = 4342.857
$ws.Range("I71").Value = 2174.625
$ws.Range("K71").Value = 10873.125
$ws.Range("M71").Value = -7129.125
$ws.Range("H82").Value = 2259.7
$ws.Range("I82").Value = 2685.4285
$ws.Range("J82").Value = 1266.3334
$ws.Range("K82").Value = 2685.4285
$ws.Range("L82").Value = 1266.3334
$ws.Range("M82").Value = -2324.4285
$ws.Range("N82").Value = -1988.3334
$ws.Range("H85").Value = 2259.7
$ws.Range("I85").Value = 2685.4285
$ws.Range("J85").Value = 1266.3334
$ws.Range("K85").Value = 2685.4285
$ws.Range("L85").Value = 1266.3334
$ws.Range("M85").Value = -1437.4285
$ws.Range("N85").Value = -3762.3334
$ws.Range("H122").Value = 7826.976
$ws.Range("I122").Value = 8294
$ws.Range("K122").Value = 24882
$ws.Range("M122").Value = -22432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9057.6
$ws.Range("I62").Value = 3994.6667
$ws.Range("J62").Value = 11227.429
$ws.Range("K62").Value = 3994.6667
$ws.Range("L62").Value = 11227.429
$ws.Range("M62").Value = -3370.6667
$ws.Range("N62").Value = -12475.429
$ws.Range("H65").Value = 9057.6
$ws.Range("I65").Value = 3994.6667
$ws.Range("J65").Value = 11227.429
$ws.Range("K65").Value = 19973.3335
$ws.Range("L65").Value = 56137.145
$ws.Range("M65").Value = -16853.3335
$ws.Range("N65").Value = -62377.145
$ws.Range("H81").Value = 4467791.5
$ws.Range("I81").Value = 6494969.5
$ws.Range("K81").Value = 12989939
$ws.Range("M81").Value = -12988878
$ws.Range("H84").Value = 4467791.5
$ws.Range("I84").Value = 6494969.5
$ws.Range("K84").Value = 64949695
$ws.Range("M84").Value = -64944391
$ws.Range("H107").Value = 5918.7334
$ws.Range("I107").Value = 1280.5555
$ws.Range("J107").Value = 12876
$ws.Range("K107").Value = 3841.6665
$ws.Range("L107").Value = 38628
$ws.Range("M107").Value = -1921.6665
$ws.Range("N107").Value = -42468
$ws.Range("H122").Value = 2878
$ws.Range("I122").Value = 2794.926
$ws.Range("K122").Value = 8384.778
$ws.Range("M122").Value = -5934.778
$ws.Range("H125").Value = 69857.10000000001
$ws.Range("J125").Value = 69857.10000000001
$ws.Range("L125").Value = 69857.10000000001
$ws.Range("N125").Value = -79697.10000000001
